$wb = $excel.ActiveWorkbook

# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" -> fix E1 label to "2050" and remove the Total row (row 13)
$sheetNames1 = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)
foreach ($name in $sheetNames1) {
    $ws = $wb.Worksheets.Item($name)
    # Leading apostrophe forces the numeric-looking label to be stored as
    # text (matching the intended "2050" text label) instead of being
    # auto-coerced into the number 2050.
    $ws.Cells.Item(1, 5).Value = "'2050"
    $ws.Rows.Item(13).Delete()
}

# Sheet 4: "Potencia Incremental - SIN(MW)" -> fix E1 label to "2041-2050" and
# remove the Total row (row 13)
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Cells.Item(1, 5).Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5 "Emissoes Totais (MtCO2eq)" is left untouched.

# Sheet 6: "Custo Total (bilhoes de R$)" -> remove the Total row (row 4)
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
